{"js": "// Load all paragraphs in the document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1. Make the \"TASKS DONE\" heading bold ---\nconst heading = paragraphs.items[0];\nheading.font.bold = true;\nheading.font.boldBidirectional = true;\n\n// --- 2. Turn the trailing (empty, bulleted) paragraph into the GitHub link line ---\nconst last = paragraphs.items[paragraphs.items.length - 1];\n\n// Remove the bullet/list numbering and list-paragraph style from this paragraph.\nlast.detachFromList();\nlast.style = \"Normal\";\n\n// Add the label text, then the hyperlink text turned into a real hyperlink.\nlast.insertText(\"LINK TO MY GITHUB WITH SOURCE CODE FOR PROJECT: \", \"End\");\nconst linkRange = last.insertText(\"https://github.com/bravine6/Sharelyft\", \"End\");\nlinkRange.hyperlink = \"https://github.com/bravine6/Sharelyft\";\n\n// --- 3. Leave a new empty paragraph behind it (matches the trailing blank line) ---\nlast.insertParagraph(\"\", \"After\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Make the \"TASKS DONE\" heading bold ---\n$heading = $d.Paragraphs.Item(1)\n$headingRange = $heading.Range\n$headingRange.Font.Bold = 1\n$headingRange.Font.BoldBi = 1\n\n# --- 2. Turn the trailing (empty, bulleted) paragraph into the GitHub link line ---\n$lastIndex = $d.Paragraphs.Count\n$last = $d.Paragraphs.Item($lastIndex)\n$lastRange = $last.Range\n\n# Remove the bullet/list numbering and list-paragraph style from this paragraph.\n$lastRange.ListFormat.RemoveNumbers()\n$last.Style = \"Normal\"\n\n# Turn the (now empty) paragraph range into a real hyperlink to the GitHub repo,\n# then prepend the label text in front of it.\n$d.Hyperlinks.Add($lastRange, \"https://github.com/bravine6/Sharelyft\", $null, $null, \"https://github.com/bravine6/Sharelyft\") | Out-Null\n$newLast = $d.Paragraphs.Item($lastIndex)\n$newLast.Range.InsertBefore(\"LINK TO MY GITHUB WITH SOURCE CODE FOR PROJECT: \")\n\n# --- 3. Leave a new empty paragraph behind it (matches the trailing blank line) ---\n$finalLast = $d.Paragraphs.Item($lastIndex)\n$finalLast.Range.InsertParagraphAfter()\n"}
